$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("2025")
$ws.Range("N2").Value = 7098.474579830174
$ws.Range("O2").Value = 6965.1427089074

$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 5707.815717280662
$ws.Range("I2").Value = 44492.05901988943
$ws.Range("L2").Value = 66334.06707325629
$ws.Range("M2").Value = 21991.42050229464
$ws.Range("N2").Value = 10503.58911838768
$ws.Range("O2").Value = 12039.81600698566

$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 2927.360317916481
$ws.Range("B2").Value = 7940.887964949257
$ws.Range("E2").Value = 67179.99183625776
$ws.Range("I2").Value = 59530.75343380851
$ws.Range("L2").Value = 66334.06707325629
$ws.Range("M2").Value = 25547.11936466757
$ws.Range("N2").Value = 15023.56613839683
$ws.Range("O2").Value = 14720.97711770872

$ws = $wb.Worksheets.Item("2040")
$ws.Range("A2").Value = 2927.360317916481
$ws.Range("B2").Value = 7940.887964949257
$ws.Range("E2").Value = 67179.99183625776
$ws.Range("I2").Value = 59530.75343380851
$ws.Range("L2").Value = 66334.06707325629
$ws.Range("M2").Value = 25547.11936466757
$ws.Range("N2").Value = 15130.35718193136
$ws.Range("O2").Value = 14720.97711770872

$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 6352.985609279765
$ws.Range("B2").Value = 7940.887964949257
$ws.Range("E2").Value = 67179.99183625776
$ws.Range("I2").Value = 59530.75343380851
$ws.Range("L2").Value = 66334.06707325629
$ws.Range("M2").Value = 25547.11936466757
$ws.Range("N2").Value = 15673.28955718277
$ws.Range("O2").Value = 17052.53181919008

$ws = $wb.Worksheets.Item("2050")
$ws.Range("A2").Value = 6352.985609279765
$ws.Range("B2").Value = 7940.887964949257
$ws.Range("E2").Value = 67179.99183625776
$ws.Range("I2").Value = 59530.75343380851
$ws.Range("L2").Value = 66334.06707325629
$ws.Range("M2").Value = 25547.11936466757
$ws.Range("N2").Value = 15673.28955718277
$ws.Range("O2").Value = 17052.53181919008

